$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restock every ingredient except Jalapenos (column I) back up to 1000.
for ($col = 1; $col -le 8; $col++) {
    $ws.Cells.Item(2, $col).Value = 1000
}
